$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previously-last row (48) had its phone number typed in as text;
# normalize it to a real number like every other row in the sheet.
$ws.Range("A48").Value = 79174445

# Append the new payment row (49): phone 71076783, Cash, 2025-08-18T17:49:44
# The phone number is entered as text (it keeps a leading-zero-safe string
# form like the other freshly-typed rows), so force a text format before
# assigning it, then drop back to the default style so no extra number
# formatting lingers on the cell.
$ws.Range("A49").NumberFormat = "@"
$ws.Range("A49").Value = "71076783"
$ws.Range("A49").Style = "Normal"

# Blank cells still need to exist in the row (matching the sheet's usual
# "blank means inline empty string" convention), so write them as an
# empty quote-prefixed text entry and then drop the style back to normal.
$ws.Range("B49").Value = "'"
$ws.Range("B49").Style = "Normal"

$ws.Range("C49").Value = "Cash"
$ws.Range("D49").Value = "2025-08-18T17:49:44"
$ws.Range("E49").Value = 200

$ws.Range("F49").Value = "'"
$ws.Range("F49").Style = "Normal"

$ws.Range("G49").Value = 200
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
